# style: remove excess vertical space between the "To" block and the
# "Subject" line (ATM / Cheque / AEPS letter templates).
#
# The template currently has five near-empty paragraphs between the last
# line of the "To," block ("...Bandhan Bank") and the "Subject: -" line:
# a run of plain spaces, a bold empty line, two tiny (6/4 half-point)
# underlined filler lines, and a final bare paragraph mark. They collapse
# to a single, plainly formatted blank paragraph (default spacing, no
# paragraph style / direct formatting).

$d = $word.ActiveDocument

# Returns the paragraph whose range starts exactly at character offset $pos.
function Get-ParagraphStartingAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -eq $pos) {
            return $p
        }
    }
    return $null
}

# Returns the paragraph whose range contains character offset $pos.
function Get-ParagraphAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $p.Range.End -ge $pos) {
            return $p
        }
    }
    return $null
}

# Locate the end of the "Bandhan Bank" paragraph (last line of the To: block).
$bandhanRng = $d.Content
$null = $bandhanRng.Find.Execute("Bandhan Bank", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bandhanPara = Get-ParagraphAt $d $bandhanRng.End
$blockStart = $bandhanPara.Range.End

# Locate the start of the "Subject: -" paragraph.
$subjectRng = $d.Content
$null = $subjectRng.Find.Execute("Subject: -", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockEnd = $subjectRng.Start - 1

# Delete everything in between except the very last paragraph mark, which
# stays in place and becomes the single blank spacer paragraph.
$r = $d.Range($blockStart, $blockEnd)
$r.Delete()

# That surviving paragraph mark keeps no style/direct formatting by default,
# but make the (already-default) spacing explicit, matching the target markup:
# no "before", no "after", single (auto) line spacing.
$remaining = Get-ParagraphStartingAt $d $blockStart
$remaining.Range.ParagraphFormat.SpaceBefore = 0
$remaining.Range.ParagraphFormat.SpaceAfter = 0
$remaining.Range.ParagraphFormat.LineSpacingRule = 0
$remaining.Range.ParagraphFormat.LineSpacing = 12
